# Fruta / hortaliza, semanal
# Update weekly price-report rows (Damasco, Feria Lagunitas de Puerto Montt)
# Rows 4, 6, 7, 8, 9 get their Fecha/Calidad/Volumen/Precio/Origen data
# rotated to reflect the latest weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("D4").Value = 44159
$ws.Range("L4").Value = "Tercera"
$ws.Range("M4").Value = 400
$ws.Range("N4").Value = 15500
$ws.Range("P4").Value = 15750
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1050

# --- Row 6 ---
$ws.Range("D6").Value = 44169
$ws.Range("M6").Value = 500
$ws.Range("R6").Value = "Región de O'Higgins"

# --- Row 7 ---
$ws.Range("D7").Value = 44176
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 500
$ws.Range("N7").Value = 15000
$ws.Range("P7").Value = 15500
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 1033

# --- Row 8 ---
$ws.Range("D8").Value = 44166
$ws.Range("M8").Value = 600
$ws.Range("N8").Value = 16000
$ws.Range("O8").Value = 17000
$ws.Range("P8").Value = 16500
$ws.Range("S8").Value = 1100

# --- Row 9 ---
$ws.Range("D9").Value = 44194
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 15500
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 1033
